$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, date range) ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Copy-format text cells (value+style already present verbatim at source) ---
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C28").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C29").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("C30").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

# --- Numeric cells that change style (text -> #,##0 or #,##0.0 number format) ---
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("E15").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E15").Value = -50
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 2
$ws.Range("H15").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H15").Value = 100
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 2
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 2
$ws.Range("E26").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E26").Value = -50
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 2
$ws.Range("H26").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H26").Value = 100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 11.111111111111
$ws.Range("L15").Value = 900
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = 42.857142857142
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = -9.090909090909
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -55.882352941176
$ws.Range("N16").Value = -89.090909090909
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 600
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 46.666666666666
$ws.Range("I17").Value = 73
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = 1.388888888888
$ws.Range("L17").Value = 92.105263157894
$ws.Range("M17").Value = 46
$ws.Range("N17").Value = -38.135593220339
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 42
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -32.258064516129
$ws.Range("L18").Value = -2.325581395348
$ws.Range("M18").Value = -70.212765957446
$ws.Range("N18").Value = -93.322734499205
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -6.382978723404
$ws.Range("I19").Value = 199
$ws.Range("J19").Value = 199
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 29.220779220779
$ws.Range("M19").Value = 41.134751773049
$ws.Range("N19").Value = -19.433198380566
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 12
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 52
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = 52.941176470588
$ws.Range("L20").Value = 79.310344827586
$ws.Range("M20").Value = -8.771929824561
$ws.Range("N20").Value = -92.905866302864
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -15.384615384615
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 9.090909090909
$ws.Range("I21").Value = 407
$ws.Range("J21").Value = 409
$ws.Range("K21").Value = -0.488997555012
$ws.Range("L21").Value = 43.816254416961
$ws.Range("M21").Value = -11.328976034858
$ws.Range("N21").Value = -79.791459781529
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -20
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 5.555555555555
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 85
$ws.Range("H24").Value = 16.470588235294
$ws.Range("I24").Value = 416
$ws.Range("J24").Value = 409
$ws.Range("K24").Value = 1.711491442542
$ws.Range("L24").Value = 69.795918367346
$ws.Range("M24").Value = 32.907348242811
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = 84.210526315789
$ws.Range("I25").Value = 119
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = 7.207207207207
$ws.Range("L25").Value = 17.821782178217
$ws.Range("M25").Value = -27.439024390243
$ws.Range("F26").Value = 4
$ws.Range("I26").Value = 11
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 266.666666666667
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 19
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -20.833333333333
$ws.Range("L27").Value = -29.629629629629
